$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style (style index 1) from an existing date cell
# so new C/D cells reuse the existing numFmtId=22 style instead of creating a new custom format.
[void]$ws.Range('C58').Copy()

# Append 47 new rows (59-105) of air-quality / toilet observation data (GT_012_W .. GT_030_M)
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = 'GT_012_W'
$ws.Cells.Item(59, 3).Value = 44771.442361111112
$ws.Cells.Item(59, 4).Value = 44771.444282407407
$ws.Range('C59:D59').PasteSpecial(-4122)
$ws.Cells.Item(59, 5).Value = 'W'
$ws.Cells.Item(59, 6).Value = 'U'

$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = 'GT_012_W'
$ws.Cells.Item(60, 3).Value = 44771.444444444445
$ws.Cells.Item(60, 4).Value = 44771.446006944447
$ws.Range('C60:D60').PasteSpecial(-4122)
$ws.Cells.Item(60, 5).Value = 'W'
$ws.Cells.Item(60, 6).Value = 'D'

$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = 'GT_012_W'
$ws.Cells.Item(61, 3).Value = 44771.447222222225
$ws.Cells.Item(61, 4).Value = 44771.449317129627
$ws.Range('C61:D61').PasteSpecial(-4122)
$ws.Cells.Item(61, 5).Value = 'W'
$ws.Cells.Item(61, 6).Value = 'M'

$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = 'GT_013_M'
$ws.Cells.Item(62, 3).Value = 44771.46597222222
$ws.Cells.Item(62, 4).Value = 44771.46707175926
$ws.Range('C62:D62').PasteSpecial(-4122)
$ws.Cells.Item(62, 5).Value = 'M'
$ws.Cells.Item(62, 6).Value = 'U'

$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = 'GT_013_M'
$ws.Cells.Item(63, 3).Value = 44771.468055555553
$ws.Cells.Item(63, 4).Value = 44771.471990740742
$ws.Range('C63:D63').PasteSpecial(-4122)
$ws.Cells.Item(63, 5).Value = 'M'
$ws.Cells.Item(63, 6).Value = 'D'

$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = 'GT_014_W'
$ws.Cells.Item(64, 3).Value = 44771.473611111112
$ws.Cells.Item(64, 4).Value = 44771.475138888891
$ws.Range('C64:D64').PasteSpecial(-4122)
$ws.Cells.Item(64, 5).Value = 'W'
$ws.Cells.Item(64, 6).Value = 'U'

$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = 'GT_014_W'
$ws.Cells.Item(65, 3).Value = 44771.476388888892
$ws.Cells.Item(65, 4).Value = 44771.482314814813
$ws.Range('C65:D65').PasteSpecial(-4122)
$ws.Cells.Item(65, 5).Value = 'W'
$ws.Cells.Item(65, 6).Value = 'D'

$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = 'GT_014_W'
$ws.Cells.Item(66, 3).Value = 44771.48333333333
$ws.Cells.Item(66, 4).Value = 44771.484780092593
$ws.Range('C66:D66').PasteSpecial(-4122)
$ws.Cells.Item(66, 5).Value = 'W'
$ws.Cells.Item(66, 6).Value = 'M'

$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = 'GT_015_W'
$ws.Cells.Item(67, 3).Value = 44771.493055555555
$ws.Cells.Item(67, 4).Value = 44771.494328703702
$ws.Range('C67:D67').PasteSpecial(-4122)
$ws.Cells.Item(67, 5).Value = 'W'
$ws.Cells.Item(67, 6).Value = 'U'

$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = 'GT_015_W'
$ws.Cells.Item(68, 3).Value = 44771.495138888888
$ws.Cells.Item(68, 4).Value = 44771.496863425928
$ws.Range('C68:D68').PasteSpecial(-4122)
$ws.Cells.Item(68, 5).Value = 'W'
$ws.Cells.Item(68, 6).Value = 'D'

$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = 'GT_015_W'
$ws.Cells.Item(69, 3).Value = 44771.497916666667
$ws.Cells.Item(69, 4).Value = 44771.499837962961
$ws.Range('C69:D69').PasteSpecial(-4122)
$ws.Cells.Item(69, 5).Value = 'W'
$ws.Cells.Item(69, 6).Value = 'M'

$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = 'GT_016_W'
$ws.Cells.Item(70, 3).Value = 44771.503472222219
$ws.Cells.Item(70, 4).Value = 44771.504340277781
$ws.Range('C70:D70').PasteSpecial(-4122)
$ws.Cells.Item(70, 5).Value = 'W'
$ws.Cells.Item(70, 6).Value = 'U'

$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = 'GT_016_W'
$ws.Cells.Item(71, 3).Value = 44771.505555555559
$ws.Cells.Item(71, 4).Value = 44771.506990740738
$ws.Range('C71:D71').PasteSpecial(-4122)
$ws.Cells.Item(71, 5).Value = 'W'
$ws.Cells.Item(71, 6).Value = 'D'

$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = 'GT_016_W'
$ws.Cells.Item(72, 3).Value = 44771.507638888892
$ws.Cells.Item(72, 4).Value = 44771.508692129632
$ws.Range('C72:D72').PasteSpecial(-4122)
$ws.Cells.Item(72, 5).Value = 'W'
$ws.Cells.Item(72, 6).Value = 'M'

$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = 'GT_017_W'
$ws.Cells.Item(73, 3).Value = 44771.57916666667
$ws.Cells.Item(73, 4).Value = 44771.580243055556
$ws.Range('C73:D73').PasteSpecial(-4122)
$ws.Cells.Item(73, 5).Value = 'W'
$ws.Cells.Item(73, 6).Value = 'U'

$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = 'GT_017_W'
$ws.Cells.Item(74, 3).Value = 44771.581250000003
$ws.Cells.Item(74, 4).Value = 44771.541030092594
$ws.Range('C74:D74').PasteSpecial(-4122)
$ws.Cells.Item(74, 5).Value = 'W'
$ws.Cells.Item(74, 6).Value = 'D'

$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = 'GT_017_W'
$ws.Cells.Item(75, 3).Value = 44771.583460648151
$ws.Cells.Item(75, 4).Value = 44771.585289351853
$ws.Range('C75:D75').PasteSpecial(-4122)
$ws.Cells.Item(75, 5).Value = 'W'
$ws.Cells.Item(75, 6).Value = 'M'

$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = 'GT_018_W'
$ws.Cells.Item(76, 3).Value = 44771.586111111108
$ws.Cells.Item(76, 4).Value = 44771.586724537039
$ws.Range('C76:D76').PasteSpecial(-4122)
$ws.Cells.Item(76, 5).Value = 'W'
$ws.Cells.Item(76, 6).Value = 'U'

$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = 'GT_018_W'
$ws.Cells.Item(77, 2).Font.Color = 255
$ws.Cells.Item(77, 3).Value = 44771.587500000001
$ws.Cells.Item(77, 4).Value = 44771.590208333335
$ws.Range('C77:D77').PasteSpecial(-4122)
$ws.Cells.Item(77, 5).Value = 'W'
$ws.Cells.Item(77, 6).Value = 'D'

$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = 'GT_018_W'
$ws.Cells.Item(78, 3).Value = 44771.59034722222
$ws.Cells.Item(78, 4).Value = 44771.591666666667
$ws.Range('C78:D78').PasteSpecial(-4122)
$ws.Cells.Item(78, 5).Value = 'W'
$ws.Cells.Item(78, 6).Value = 'M'

$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = 'GT_019_M'
$ws.Cells.Item(79, 2).Font.Color = 0
$ws.Cells.Item(79, 3).Value = 44771.593055555553
$ws.Cells.Item(79, 4).Value = 44771.593854166669
$ws.Range('C79:D79').PasteSpecial(-4122)
$ws.Cells.Item(79, 5).Value = 'M'
$ws.Cells.Item(79, 6).Value = 'U'

$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = 'GT_019_M'
$ws.Cells.Item(80, 2).Font.Color = 0
$ws.Cells.Item(80, 3).Value = 44771.594444444447
$ws.Cells.Item(80, 4).Value = 44771.595520833333
$ws.Range('C80:D80').PasteSpecial(-4122)
$ws.Cells.Item(80, 5).Value = 'M'
$ws.Cells.Item(80, 6).Value = 'D'

$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = 'GT_020_M'
$ws.Cells.Item(81, 2).Font.Color = 0
$ws.Cells.Item(81, 3).Value = 44771.634722222225
$ws.Cells.Item(81, 4).Value = 44771.635949074072
$ws.Range('C81:D81').PasteSpecial(-4122)
$ws.Cells.Item(81, 5).Value = 'M'
$ws.Cells.Item(81, 6).Value = 'U'

$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = 'GT_020_M'
$ws.Cells.Item(82, 2).Font.Color = 0
$ws.Cells.Item(82, 3).Value = 44771.636805555558
$ws.Cells.Item(82, 4).Value = 44771.639907407407
$ws.Range('C82:D82').PasteSpecial(-4122)
$ws.Cells.Item(82, 5).Value = 'M'
$ws.Cells.Item(82, 6).Value = 'D'

$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = 'GT_021_M'
$ws.Cells.Item(83, 2).Font.Color = 0
$ws.Cells.Item(83, 3).Value = 44771.64166666667
$ws.Cells.Item(83, 4).Value = 44771.642361111109
$ws.Range('C83:D83').PasteSpecial(-4122)
$ws.Cells.Item(83, 5).Value = 'M'
$ws.Cells.Item(83, 6).Value = 'U'

$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = 'GT_021_M'
$ws.Cells.Item(84, 2).Font.Color = 0
$ws.Cells.Item(84, 3).Value = 44771.643055555556
$ws.Cells.Item(84, 4).Value = 44771.645312499997
$ws.Range('C84:D84').PasteSpecial(-4122)
$ws.Cells.Item(84, 5).Value = 'M'
$ws.Cells.Item(84, 6).Value = 'D'

$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = 'GT_022_M'
$ws.Cells.Item(85, 2).Font.Color = 0
$ws.Cells.Item(85, 3).Value = 44771.663888888892
$ws.Cells.Item(85, 4).Value = 44771.665277777778
$ws.Range('C85:D85').PasteSpecial(-4122)
$ws.Cells.Item(85, 5).Value = 'M'
$ws.Cells.Item(85, 6).Value = 'U'

$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = 'GT_022_M'
$ws.Cells.Item(86, 2).Font.Color = 0
$ws.Cells.Item(86, 3).Value = 44771.665972222225
$ws.Cells.Item(86, 4).Value = 44771.668923611112
$ws.Range('C86:D86').PasteSpecial(-4122)
$ws.Cells.Item(86, 5).Value = 'M'
$ws.Cells.Item(86, 6).Value = 'D'

$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = 'GT_023_W'
$ws.Cells.Item(87, 2).Font.Color = 0
$ws.Cells.Item(87, 3).Value = 44771.68472222222
$ws.Cells.Item(87, 4).Value = 44771.685763888891
$ws.Range('C87:D87').PasteSpecial(-4122)
$ws.Cells.Item(87, 5).Value = 'W'
$ws.Cells.Item(87, 6).Value = 'U'

$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = 'GT_023_W'
$ws.Cells.Item(88, 2).Font.Color = 0
$ws.Cells.Item(88, 3).Value = 44771.686805555553
$ws.Cells.Item(88, 4).Value = 44771.68854166667
$ws.Range('C88:D88').PasteSpecial(-4122)
$ws.Cells.Item(88, 5).Value = 'W'
$ws.Cells.Item(88, 6).Value = 'D'

$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = 'GT_023_W'
$ws.Cells.Item(89, 2).Font.Color = 0
$ws.Cells.Item(89, 3).Value = 44771.69027777778
$ws.Cells.Item(89, 4).Value = 44771.691944444443
$ws.Range('C89:D89').PasteSpecial(-4122)
$ws.Cells.Item(89, 5).Value = 'W'
$ws.Cells.Item(89, 6).Value = 'M'

$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = 'GT_024_M'
$ws.Cells.Item(90, 2).Font.Color = 0
$ws.Cells.Item(90, 3).Value = 44774.461111111108
$ws.Cells.Item(90, 4).Value = 44774.462337962963
$ws.Range('C90:D90').PasteSpecial(-4122)
$ws.Cells.Item(90, 5).Value = 'M'
$ws.Cells.Item(90, 6).Value = 'U'

$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = 'GT_024_M'
$ws.Cells.Item(91, 2).Font.Color = 0
$ws.Cells.Item(91, 3).Value = 44774.463194444441
$ws.Cells.Item(91, 4).Value = 44774.465578703705
$ws.Range('C91:D91').PasteSpecial(-4122)
$ws.Cells.Item(91, 5).Value = 'M'
$ws.Cells.Item(91, 6).Value = 'D'

$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = 'GT_025_W'
$ws.Cells.Item(92, 2).Font.Color = 0
$ws.Cells.Item(92, 3).Value = 44774.493055555555
$ws.Cells.Item(92, 4).Value = 44774.494467592594
$ws.Range('C92:D92').PasteSpecial(-4122)
$ws.Cells.Item(92, 5).Value = 'W'
$ws.Cells.Item(92, 6).Value = 'U'

$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = 'GT_025_W'
$ws.Cells.Item(93, 2).Font.Color = 0
$ws.Cells.Item(93, 3).Value = 44774.495833333334
$ws.Cells.Item(93, 4).Value = 44774.498472222222
$ws.Range('C93:D93').PasteSpecial(-4122)
$ws.Cells.Item(93, 5).Value = 'W'
$ws.Cells.Item(93, 6).Value = 'D'

$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = 'GT_025_W'
$ws.Cells.Item(94, 2).Font.Color = 0
$ws.Cells.Item(94, 3).Value = 44774.499305555553
$ws.Cells.Item(94, 4).Value = 44774.501840277779
$ws.Range('C94:D94').PasteSpecial(-4122)
$ws.Cells.Item(94, 5).Value = 'W'
$ws.Cells.Item(94, 6).Value = 'M'

$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = 'GT_026_W'
$ws.Cells.Item(95, 2).Font.Color = 0
$ws.Cells.Item(95, 3).Value = 44774.50277777778
$ws.Cells.Item(95, 4).Value = 44774.504293981481
$ws.Range('C95:D95').PasteSpecial(-4122)
$ws.Cells.Item(95, 5).Value = 'W'
$ws.Cells.Item(95, 6).Value = 'U'

$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = 'GT_026_W'
$ws.Cells.Item(96, 2).Font.Color = 0
$ws.Cells.Item(96, 3).Value = 44774.504861111112
$ws.Cells.Item(96, 4).Value = 44774.507141203707
$ws.Range('C96:D96').PasteSpecial(-4122)
$ws.Cells.Item(96, 5).Value = 'W'
$ws.Cells.Item(96, 6).Value = 'D'

$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = 'GT_026_W'
$ws.Cells.Item(97, 2).Font.Color = 0
$ws.Cells.Item(97, 3).Value = 44774.508333333331
$ws.Cells.Item(97, 4).Value = 44774.51059027778
$ws.Range('C97:D97').PasteSpecial(-4122)
$ws.Cells.Item(97, 5).Value = 'W'
$ws.Cells.Item(97, 6).Value = 'M'

$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = 'GT_027_M'
$ws.Cells.Item(98, 2).Font.Color = 0
$ws.Cells.Item(98, 3).Value = 44774.529861111114
$ws.Cells.Item(98, 4).Value = 44774.530729166669
$ws.Range('C98:D98').PasteSpecial(-4122)
$ws.Cells.Item(98, 5).Value = 'M'
$ws.Cells.Item(98, 6).Value = 'U'

$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = 'GT_027_M'
$ws.Cells.Item(99, 2).Font.Color = 0
$ws.Cells.Item(99, 3).Value = 44774.53125
$ws.Cells.Item(99, 4).Value = 44774.532650462963
$ws.Range('C99:D99').PasteSpecial(-4122)
$ws.Cells.Item(99, 5).Value = 'M'
$ws.Cells.Item(99, 6).Value = 'D'

$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = 'GT_028_M'
$ws.Cells.Item(100, 2).Font.Color = 0
$ws.Cells.Item(100, 3).Value = 44774.585416666669
$ws.Cells.Item(100, 4).Value = 44774.586840277778
$ws.Range('C100:D100').PasteSpecial(-4122)
$ws.Cells.Item(100, 5).Value = 'M'
$ws.Cells.Item(100, 6).Value = 'U'

$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = 'GT_028_M'
$ws.Cells.Item(101, 2).Font.Color = 0
$ws.Cells.Item(101, 3).Value = 44774.587500000001
$ws.Cells.Item(101, 4).Value = 44774.591157407405
$ws.Range('C101:D101').PasteSpecial(-4122)
$ws.Cells.Item(101, 5).Value = 'M'
$ws.Cells.Item(101, 6).Value = 'D'

$ws.Cells.Item(102, 1).Value = 101
$ws.Cells.Item(102, 2).Value = 'GT_029_M'
$ws.Cells.Item(102, 2).Font.Color = 0
$ws.Cells.Item(102, 3).Value = 44774.630555555559
$ws.Cells.Item(102, 4).Value = 44774.631944444445
$ws.Range('C102:D102').PasteSpecial(-4122)
$ws.Cells.Item(102, 5).Value = 'M'
$ws.Cells.Item(102, 6).Value = 'U'

$ws.Cells.Item(103, 1).Value = 102
$ws.Cells.Item(103, 2).Value = 'GT_029_M'
$ws.Cells.Item(103, 2).Font.Color = 0
$ws.Cells.Item(103, 3).Value = 44774.632638888892
$ws.Cells.Item(103, 4).Value = 44774.635787037034
$ws.Range('C103:D103').PasteSpecial(-4122)
$ws.Cells.Item(103, 5).Value = 'M'
$ws.Cells.Item(103, 6).Value = 'D'

$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 2).Value = 'GT_030_M'
$ws.Cells.Item(104, 2).Font.Color = 0
$ws.Cells.Item(104, 3).Value = 44774.647916666669
$ws.Cells.Item(104, 4).Value = 44774.6487037037
$ws.Range('C104:D104').PasteSpecial(-4122)
$ws.Cells.Item(104, 5).Value = 'M'
$ws.Cells.Item(104, 6).Value = 'U'

$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = 'GT_030_M'
$ws.Cells.Item(105, 2).Font.Color = 0
$ws.Cells.Item(105, 3).Value = 44774.649305555555
$ws.Cells.Item(105, 4).Value = 44774.65079861111
$ws.Range('C105:D105').PasteSpecial(-4122)
$ws.Cells.Item(105, 5).Value = 'M'
$ws.Cells.Item(105, 6).Value = 'D'

$excel.CutCopyMode = $false

# Update the view: scroll/select to show the newly added rows
[void]$ws.Range('A57:A105').Select()